$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Student Summary"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Student Summary")

# Insert the two new "Course Code" / "Max marks" rows (A11:C12), copying the
# look (style) of the row above (A9:C10) and then filling in the new text.
$ws1.Range("A9:C10").Copy()
$ws1.Range("A11:C12").PasteSpecial(-4122)

$ws1.Range("B11").Value = "Course Code:"
$ws1.Range("C11").Value = "DSPC601"
$ws1.Range("B12").Value = "Max marks"
$ws1.Range("C12").Value = 40

# Relabel the summary attribute names.
$ws1.Range("A17").Value = "Average Marks"
$ws1.Range("A18").Value = "Less Than 40%"
$ws1.Range("A19").Value = "Between 40 % - 75 %"
$ws1.Range("A20").Value = "More than 75%"

# Round the average-marks value.
$ws1.Range("B17").Value = 26.84

# ---------------------------------------------------------------------------
# Sheet 3: "Fast Learners" -- reorder existing rows and append new students
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Fast Learners")

$fastLearners = @(
    @(1, 2136110001, "Aravind S", 40),
    @(46, 2136110047, "Gowtham R", 39),
    @(29, 2136110030, "Brijesh A", 39),
    @(24, 2136110025, "Adhithi K", 39),
    @(16, 2136110016, "Nithya Sri R", 38),
    @(47, 2136110048, "Bhuvanadurai M", 38),
    @(30, 2136110031, "Hitesh Kumar K A", 38),
    @(22, 2136110023, "Surya Prakash R", 37),
    @(48, 2136110049, "Kailashwaran R", 35),
    @(49, 2236150001, "Dhanush B", 35),
    @(10, 2136110010, "Krishnapriya K", 35),
    @(4, 2136110004, "Deepakragavan J", 33),
    @(7, 2136110007, "Guruprasath V", 33),
    @(9, 2136110009, "Kalaivani S", 33),
    @(36, 2136110037, "Mohamed Tharif B", 32),
    @(37, 2136110038, "Pradeep M", 32),
    @(40, 2136110041, "Ragothaman R", 32),
    @(21, 2136110022, "Suji Shri B", 31),
    @(33, 2136110034, "Kaviraj M", 31),
    @(39, 2136110040, "Preethiga S", 31),
    @(41, 2136110042, "Rajadurai P", 30),
    @(17, 2136110018, "Sikanthkumar C", 30),
    @(18, 2136110019, "Sivaa Ganesh S", 30),
    @(44, 2136110045, "Varsha V", 29),
    @(3, 2136110003, "Ashik Jenly V L", 27),
    @(8, 2136110008, "Jananika B", 27),
    @(13, 2136110013, "Naveena A", 27),
    @(31, 2136110032, "Jaikrishnan V", 26),
    @(15, 2136110015, "Nilavanan S.A", 26),
    @(25, 2136110026, "Ajay S", 26)
)

$r = 2
foreach ($rec in $fastLearners) {
    $ws3.Cells.Item($r, 1).Value = $rec[0]
    $ws3.Cells.Item($r, 2).Value = $rec[1]
    $ws3.Cells.Item($r, 3).Value = $rec[2]
    $ws3.Cells.Item($r, 4).Value = $rec[3]
    $r = $r + 1
}
